$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 45952
$ws.Range("B12").Value = 586
$ws.Range("C12").Value = 29
$ws.Range("D12").Value = 557

$ws.Range("A12:D12").Select()
